# Fruta / hortaliza, semanal
# Insert a new weekly record (row 55) for "Choclo - Dulce o Americano"
# from Región de Arica y Parinacota, pushing the existing rows 55-57 down
# to rows 56-58.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 55; this shifts old rows 55-57 to 56-58
# and copies the formatting of the surrounding rows (e.g. the date style
# on column D) into the freshly inserted row.
$ws.Rows.Item(55).Insert()

# Populate the newly inserted row 55 with the new weekly data point.
$ws.Range("A55").Value2 = 11
$ws.Range("B55").Value2 = "Vega Monumental Concepción"
$ws.Range("C55").Value2 = "Bíobío"
$ws.Range("D55").Value2 = 44509
$ws.Range("E55").Value2 = 8
$ws.Range("F55").Value2 = 100112024
$ws.Range("G55").Value2 = "Choclo"
$ws.Range("H55").Value2 = "Dulce o Americano"
$ws.Range("I55").Value2 = "Primera"
$ws.Range("J55").Value2 = 100
$ws.Range("K55").Value2 = 20000
$ws.Range("L55").Value2 = 22000
$ws.Range("M55").Value2 = 21000
$ws.Range("N55").Value2 = "$/malla 70 unidades"
$ws.Range("O55").Value2 = "Región de Arica y Parinacota"
$ws.Range("P55").Value2 = 300
$ws.Range("Q55").Value2 = 70
$ws.Range("R55").Value2 = "Hortaliza"
